$wb = $excel.ActiveWorkbook

# --- Birds sheet: reassign birds previously pointing at the (invalid/duplicate)
# cage "a1" to the real cage "a2" ---
$birds = $wb.Worksheets.Item("Birds")
$rowsToFix = @(2, 3, 4, 7, 8, 10, 16)
foreach ($r in $rowsToFix) {
    $birds.Cells.Item($r, 8).Value = "a2"
}
[void]$birds.Range("A2").Select()

# --- Cages sheet: add the three new cages that were missing / newly created ---
$cages = $wb.Worksheets.Item("Cages")
$cages.Cells.Item(11, 1).Value = "a11"
$cages.Cells.Item(11, 2).Value = 3
$cages.Cells.Item(11, 3).Value = 3
$cages.Cells.Item(11, 4).Value = 3
$cages.Cells.Item(11, 5).Value = "plastic"

$cages.Cells.Item(12, 1).Value = "a12"
$cages.Cells.Item(12, 2).Value = -4
$cages.Cells.Item(12, 3).Value = -5
$cages.Cells.Item(12, 4).Value = -6
$cages.Cells.Item(12, 5).Value = "wood"

$cages.Cells.Item(13, 1).Value = "a13"
$cages.Cells.Item(13, 2).Value = 2
$cages.Cells.Item(13, 3).Value = 2
$cages.Cells.Item(13, 4).Value = 2
$cages.Cells.Item(13, 5).Value = "wood"
